# Updates cryptos list values (Price / Volume(1h) columns) and a few
# coin name/link swaps, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''30.496.32'
$ws.Cells.Item(2, 5).Value = '  +1.89%  '

$ws.Cells.Item(3, 4).Value = '''1.673.51'
$ws.Cells.Item(3, 5).Value = '  +2.51%  '

$ws.Cells.Item(4, 5).Value = '  -0.04%  '

$ws.Cells.Item(5, 4).Value = '''219.64'
$ws.Cells.Item(5, 5).Value = '  +2.46%  '

$ws.Cells.Item(6, 5).Value = '  +1.40%  '

$ws.Cells.Item(7, 5).Value = '  -0.03%  '

$ws.Cells.Item(8, 4).Value = '''29.70'
$ws.Cells.Item(8, 5).Value = '  +4.17%  '

$ws.Cells.Item(9, 4).Value = '''0.264'
$ws.Cells.Item(9, 5).Value = '  +2.43%  '

$ws.Cells.Item(10, 4).Value = '''0.0636'
$ws.Cells.Item(10, 5).Value = '  +4.69%  '

$ws.Cells.Item(11, 4).Value = '''0.0904'
$ws.Cells.Item(11, 5).Value = '  -0.91%  '

$ws.Cells.Item(12, 4).Value = '''1.912.07'
$ws.Cells.Item(12, 5).Value = '  +2.45%  '

$ws.Cells.Item(13, 4).Value = '''1.677.03'
$ws.Cells.Item(13, 5).Value = '  +2.81%  '

$ws.Cells.Item(14, 4).Value = '''0.612'
$ws.Cells.Item(14, 5).Value = '  +8.84%  '

$ws.Cells.Item(15, 4).Value = '''10.09'
$ws.Cells.Item(15, 5).Value = '  +8.95%  '

$ws.Cells.Item(16, 4).Value = '''3.97'
$ws.Cells.Item(16, 5).Value = '  +3.18%  '

$ws.Cells.Item(17, 4).Value = '''30.520.88'
$ws.Cells.Item(17, 5).Value = '  +1.94%  '

$ws.Cells.Item(18, 4).Value = '''66.12'
$ws.Cells.Item(18, 5).Value = '  +3.38%  '

$ws.Cells.Item(19, 4).Value = '''242.39'
$ws.Cells.Item(19, 5).Value = '  -0.03%  '

$ws.Cells.Item(20, 4).Value = '''0.0₃0719'
$ws.Cells.Item(20, 5).Value = '  +2.55%  '

$ws.Cells.Item(21, 5).Value = '  +0.01%  '

$ws.Cells.Item(22, 4).Value = '''4.25'
$ws.Cells.Item(22, 5).Value = '  +2.94%  '

$ws.Cells.Item(23, 4).Value = '''9.95'
$ws.Cells.Item(23, 5).Value = '  +0.98%  '

$ws.Cells.Item(24, 4).Value = '''2.14'
$ws.Cells.Item(24, 5).Value = '  -0.39%  '

$ws.Cells.Item(25, 4).Value = '''158.16'
$ws.Cells.Item(25, 5).Value = '  +0.35%  '

$ws.Cells.Item(26, 4).Value = '''15.83'
$ws.Cells.Item(26, 5).Value = '  +2.11%  '

$ws.Cells.Item(27, 5).Value = '  +2.13%  '

$ws.Cells.Item(28, 4).Value = '''6.67'
$ws.Cells.Item(28, 5).Value = '  +1.06%  '

$ws.Cells.Item(29, 5).Value = '  -0.12%  '

$ws.Cells.Item(30, 4).Value = '''0.0493'
$ws.Cells.Item(30, 5).Value = '  +1.53%  '

$ws.Cells.Item(31, 5).Value = '  +2.85%  '

$ws.Cells.Item(32, 5).Value = '  +2.69%  '

$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 4).Value = '''3.29'
$ws.Cells.Item(33, 5).Value = '  +3.60%  '

$ws.Cells.Item(34, 2).Value = 'Maker'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(34, 4).Value = '''1.500.40'
$ws.Cells.Item(34, 5).Value = '  +5.41%  '

$ws.Cells.Item(35, 4).Value = '''1.76'
$ws.Cells.Item(35, 5).Value = '  +6.99%  '

$ws.Cells.Item(36, 4).Value = '''84.25'
$ws.Cells.Item(36, 5).Value = '  +11.20%  '

$ws.Cells.Item(37, 5).Value = '  -1.17%  '

$ws.Cells.Item(38, 2).Value = 'ImmutableX'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(38, 4).Value = '''0.595'
$ws.Cells.Item(38, 5).Value = '  +7.84%  '

$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).Value = '''0.0178'
$ws.Cells.Item(39, 5).Value = '  +5.33%  '

$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(40, 4).Value = '''2.67'
$ws.Cells.Item(40, 5).Value = '  -4.44%  '

$ws.Cells.Item(41, 2).Value = 'HuobiToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(41, 4).Value = '''2.29'
$ws.Cells.Item(41, 5).Value = '  -0.22%  '

$ws.Cells.Item(42, 4).Value = '''0.838'
$ws.Cells.Item(42, 5).Value = '  +1.45%  '

$ws.Cells.Item(43, 5).Value = '  +1.80%  '

$ws.Cells.Item(44, 5).Value = '  -1.21%  '

$ws.Cells.Item(45, 5).Value = '  +0.14%  '

$ws.Cells.Item(46, 5).Value = '  -0.02%  '

$ws.Cells.Item(47, 4).Value = '''5.53'
$ws.Cells.Item(47, 5).Value = '  +3.41%  '

$ws.Cells.Item(48, 4).Value = '''51.12'
$ws.Cells.Item(48, 5).Value = '  -3.68%  '

$ws.Cells.Item(49, 4).Value = '''1.802.43'
$ws.Cells.Item(49, 5).Value = '  +1.60%  '

$ws.Cells.Item(50, 4).Value = '''94.59'
$ws.Cells.Item(50, 5).Value = '  +5.55%  '

$ws.Cells.Item(51, 4).Value = '''0.0₆0113'
$ws.Cells.Item(51, 5).Value = '  -0.03%  '
